$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing data row (row 2) into a new row 3 before changing
# row 2's browser value, so row 3 keeps the original "chrome" entry.
$ws.Range("A3").Value = $ws.Range("A2").Value()
$ws.Range("B3").Value = $ws.Range("B2").Value()
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Range("D3").Value = $ws.Range("D2").Value()

# Copy the formatting of row 2 down to row 3 so styles match.
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)

# Update row 2's browser to firefox.
$ws.Range("D2").Value = "firefox"

# Restore the selection so it matches the saved worksheet view.
$ws.Range("F10").Select() | Out-Null
